# Auto-generated script to update Typhon_Profits workbook values
# per commit: chore: update Sheets via scheduled runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 13832
$ws.Cells.Item(17, 10).Value = 14749.923
$ws.Cells.Item(17, 12).Value = 44249.769
$ws.Cells.Item(17, 14).Value = -44585.769
$ws.Cells.Item(112, 8).Value = 3087564
$ws.Cells.Item(112, 10).Value = 3087564
$ws.Cells.Item(112, 12).Value = 9262692
$ws.Cells.Item(112, 14).Value = -9264908
$ws.Cells.Item(113, 8).Value = 66671900
$ws.Cells.Item(113, 9).Value = 90913140
$ws.Cells.Item(113, 10).Value = 8500
$ws.Cells.Item(113, 11).Value = 90913140
$ws.Cells.Item(113, 12).Value = 8500
$ws.Cells.Item(113, 13).Value = -90909886
$ws.Cells.Item(113, 14).Value = -15008
$ws.Cells.Item(116, 8).Value = 3995.2307
$ws.Cells.Item(116, 9).Value = 2074.4443
$ws.Cells.Item(116, 10).Value = 5012.1177
$ws.Cells.Item(116, 11).Value = 2074.4443
$ws.Cells.Item(116, 12).Value = 5012.1177
$ws.Cells.Item(116, 13).Value = 1367.5557
$ws.Cells.Item(116, 14).Value = -11896.1177
$ws.Cells.Item(129, 8).Value = 1082.9574
$ws.Cells.Item(129, 9).Value = 373.33334
$ws.Cells.Item(129, 10).Value = 1251.0264
$ws.Cells.Item(129, 11).Value = 1120.00002
$ws.Cells.Item(129, 12).Value = 3753.0792
$ws.Cells.Item(129, 13).Value = 3879.99998
$ws.Cells.Item(129, 14).Value = -13753.0792
$ws.Cells.Item(132, 8).Value = 3931.5
$ws.Cells.Item(132, 9).Value = 4495.684
$ws.Cells.Item(132, 10).Value = 2400.1428
$ws.Cells.Item(132, 11).Value = 13487.052
$ws.Cells.Item(132, 12).Value = 7200.428400000001
$ws.Cells.Item(132, 13).Value = -10957.052
$ws.Cells.Item(132, 14).Value = -12260.4284
$ws.Cells.Item(138, 8).Value = 33336664
$ws.Cells.Item(138, 9).Value = 66668844
$ws.Cells.Item(138, 10).Value = 4483.8
$ws.Cells.Item(138, 11).Value = 200006532
$ws.Cells.Item(138, 12).Value = 13451.4
$ws.Cells.Item(138, 13).Value = -200001392
$ws.Cells.Item(138, 14).Value = -23731.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4184.477
$ws.Cells.Item(32, 9).Value = 3287.7627
$ws.Cells.Item(32, 11).Value = 3287.7627
$ws.Cells.Item(32, 13).Value = -3000.7627
$ws.Cells.Item(102, 8).Value = 700.9
$ws.Cells.Item(102, 9).Value = 658.8946999999999
$ws.Cells.Item(102, 10).Value = 1499
$ws.Cells.Item(102, 11).Value = 658.8946999999999
$ws.Cells.Item(102, 12).Value = 1499
$ws.Cells.Item(102, 13).Value = 963.1053000000001
$ws.Cells.Item(102, 14).Value = -4743
$ws.Cells.Item(122, 8).Value = 2577.3704
$ws.Cells.Item(122, 9).Value = 2198.7058
$ws.Cells.Item(122, 10).Value = 3221.1
$ws.Cells.Item(122, 11).Value = 6596.117400000001
$ws.Cells.Item(122, 12).Value = 9663.299999999999
$ws.Cells.Item(122, 13).Value = -4146.117400000001
$ws.Cells.Item(122, 14).Value = -14563.3
$ws.Cells.Item(132, 8).Value = 13100.023
$ws.Cells.Item(132, 9).Value = 1139.0625
$ws.Cells.Item(132, 10).Value = 47895.547
$ws.Cells.Item(132, 11).Value = 3417.1875
$ws.Cells.Item(132, 12).Value = 143686.641
$ws.Cells.Item(132, 13).Value = -887.1875
$ws.Cells.Item(132, 14).Value = -148746.641

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2533.3333
$ws.Cells.Item(20, 9).Value = 2550
$ws.Cells.Item(20, 10).Value = 2500
$ws.Cells.Item(20, 11).Value = 2550
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = -2303
$ws.Cells.Item(20, 14).Value = -2994
$ws.Cells.Item(134, 8).Value = 3803.3872
$ws.Cells.Item(134, 9).Value = 4032.3215
$ws.Cells.Item(134, 11).Value = 12096.9645
$ws.Cells.Item(134, 13).Value = -9561.9645

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 753.9286
$ws.Cells.Item(16, 9).Value = 696.7
$ws.Cells.Item(16, 10).Value = 897
$ws.Cells.Item(16, 11).Value = 696.7
$ws.Cells.Item(16, 12).Value = 897
$ws.Cells.Item(16, 13).Value = -409.7
$ws.Cells.Item(16, 14).Value = -1471
$ws.Cells.Item(22, 8).Value = 372.75
$ws.Cells.Item(22, 9).Value = 280
$ws.Cells.Item(22, 10).Value = 651
$ws.Cells.Item(22, 11).Value = 280
$ws.Cells.Item(22, 12).Value = 651
$ws.Cells.Item(22, 13).Value = 70
$ws.Cells.Item(22, 14).Value = -1351
$ws.Cells.Item(99, 8).Value = 25003540
$ws.Cells.Item(99, 9).Value = 2968.7693
$ws.Cells.Item(99, 10).Value = 71433170
$ws.Cells.Item(99, 11).Value = 2968.7693
$ws.Cells.Item(99, 12).Value = 71433170
$ws.Cells.Item(99, 13).Value = -1470.7693
$ws.Cells.Item(99, 14).Value = -71436166
$ws.Cells.Item(107, 8).Value = 1154.2727
$ws.Cells.Item(107, 9).Value = 837.8
$ws.Cells.Item(107, 11).Value = 837.8
$ws.Cells.Item(107, 13).Value = 1082.2
$ws.Cells.Item(113, 8).Value = 753.9286
$ws.Cells.Item(113, 9).Value = 696.7
$ws.Cells.Item(113, 10).Value = 897
$ws.Cells.Item(113, 11).Value = 696.7
$ws.Cells.Item(113, 12).Value = 897
$ws.Cells.Item(113, 13).Value = 1473.3
$ws.Cells.Item(113, 14).Value = -5237
$ws.Cells.Item(126, 8).Value = 25003540
$ws.Cells.Item(126, 9).Value = 2968.7693
$ws.Cells.Item(126, 10).Value = 71433170
$ws.Cells.Item(126, 11).Value = 8906.3079
$ws.Cells.Item(126, 12).Value = 214299510
$ws.Cells.Item(126, 13).Value = -6436.3079
$ws.Cells.Item(126, 14).Value = -214304450

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 463.64285
$ws.Cells.Item(92, 9).Value = 266.8889
$ws.Cells.Item(92, 10).Value = 817.8
$ws.Cells.Item(92, 11).Value = 800.6667
$ws.Cells.Item(92, 12).Value = 2453.4
$ws.Cells.Item(92, 13).Value = 447.3333
$ws.Cells.Item(92, 14).Value = -4949.4
$ws.Cells.Item(113, 8).Value = 389.22223
$ws.Cells.Item(113, 9).Value = 387.18182
$ws.Cells.Item(113, 10).Value = 392.42856
$ws.Cells.Item(113, 11).Value = 1161.54546
$ws.Cells.Item(113, 12).Value = 1177.28568
$ws.Cells.Item(113, 13).Value = 1008.45454
$ws.Cells.Item(113, 14).Value = -5517.28568
$ws.Cells.Item(123, 8).Value = 4410
$ws.Cells.Item(123, 9).Value = 980
$ws.Cells.Item(123, 10).Value = 5267.5
$ws.Cells.Item(123, 11).Value = 2940
$ws.Cells.Item(123, 12).Value = 15802.5
$ws.Cells.Item(123, 13).Value = -490
$ws.Cells.Item(123, 14).Value = -20702.5
$ws.Cells.Item(131, 8).Value = 724.42
$ws.Cells.Item(131, 9).Value = 376.55554
$ws.Cells.Item(131, 10).Value = 758.82416
$ws.Cells.Item(131, 11).Value = 1129.66662
$ws.Cells.Item(131, 12).Value = 2276.47248
$ws.Cells.Item(131, 13).Value = 3910.33338
$ws.Cells.Item(131, 14).Value = -12356.47248

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8945857
$ws.Cells.Item(70, 10).Value = 8945857
$ws.Cells.Item(70, 12).Value = 8945857
$ws.Cells.Item(70, 14).Value = -8946397
$ws.Cells.Item(73, 8).Value = 8945857
$ws.Cells.Item(73, 10).Value = 8945857
$ws.Cells.Item(73, 12).Value = 8945857
$ws.Cells.Item(73, 14).Value = -8947729
$ws.Cells.Item(102, 8).Value = 20002632
$ws.Cells.Item(102, 9).Value = 22729696
$ws.Cells.Item(102, 10).Value = 4171.3335
$ws.Cells.Item(102, 11).Value = 22729696
$ws.Cells.Item(102, 12).Value = 4171.3335
$ws.Cells.Item(102, 13).Value = -22728074
$ws.Cells.Item(102, 14).Value = -7415.3335
$ws.Cells.Item(113, 8).Value = 2280.1904
$ws.Cells.Item(113, 9).Value = 1819.8462
$ws.Cells.Item(113, 10).Value = 3028.25
$ws.Cells.Item(113, 11).Value = 1819.8462
$ws.Cells.Item(113, 12).Value = 3028.25
$ws.Cells.Item(113, 13).Value = 350.1538
$ws.Cells.Item(113, 14).Value = -7368.25
$ws.Cells.Item(122, 8).Value = 102566040
$ws.Cells.Item(122, 9).Value = 41668816
$ws.Cells.Item(122, 10).Value = 200001600
$ws.Cells.Item(122, 11).Value = 125006448
$ws.Cells.Item(122, 12).Value = 600004800
$ws.Cells.Item(122, 13).Value = -125003998
$ws.Cells.Item(122, 14).Value = -600009700
$ws.Cells.Item(126, 8).Value = 4282.4863
$ws.Cells.Item(126, 9).Value = 3561.158
$ws.Cells.Item(126, 10).Value = 5043.8887
$ws.Cells.Item(126, 11).Value = 10683.474
$ws.Cells.Item(126, 12).Value = 15131.6661
$ws.Cells.Item(126, 13).Value = -8213.474
$ws.Cells.Item(126, 14).Value = -20071.6661
$ws.Cells.Item(132, 8).Value = 18427.867
$ws.Cells.Item(132, 9).Value = 1291.619
$ws.Cells.Item(132, 10).Value = 58412.445
$ws.Cells.Item(132, 11).Value = 3874.857
$ws.Cells.Item(132, 12).Value = 175237.335
$ws.Cells.Item(132, 13).Value = -1344.857
$ws.Cells.Item(132, 14).Value = -180297.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 702814.25
$ws.Cells.Item(122, 9).Value = 893463.6
$ws.Cells.Item(122, 10).Value = 3766.5
$ws.Cells.Item(122, 11).Value = 2680390.8
$ws.Cells.Item(122, 12).Value = 11299.5
$ws.Cells.Item(122, 13).Value = -2677940.8
$ws.Cells.Item(122, 14).Value = -16199.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1014.68
$ws.Cells.Item(122, 9).Value = 989.4545000000001
$ws.Cells.Item(122, 11).Value = 2968.3635
$ws.Cells.Item(122, 13).Value = -518.3635000000004
$ws.Cells.Item(132, 8).Value = 1096.5385
$ws.Cells.Item(132, 9).Value = 743.05554
$ws.Cells.Item(132, 10).Value = 1891.875
$ws.Cells.Item(132, 11).Value = 2229.16662
$ws.Cells.Item(132, 12).Value = 5675.625
$ws.Cells.Item(132, 13).Value = 300.83338
$ws.Cells.Item(132, 14).Value = -10735.625
